# Update prevalence2018 sheet: new STEP-data-aligned prevalence values for
# ages 19-120 (rows 21-122, column C). The incidence2018_plus sheet pulls
# these via formulas (=prevalence2018!Cxx/100) so it recalculates itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prevalence2018")

# age 19-24 (rows 21-26): 0.35 -> 0
$ws.Range("C21:C26").Value = 0

# age 25-34 (rows 27-36): 0.35 -> 0.36
$ws.Range("C27:C36").Value = 0.36

# age 35-44 (rows 37-46): 0.43 -> 0.425
$ws.Range("C37:C46").Value = 0.425

# age 45-54 (rows 47-56): 0.57 -> 0.5
$ws.Range("C47:C56").Value = 0.5

# age 55-120 (rows 57-122): 0.9 -> 0.97
$ws.Range("C57:C122").Value = 0.97

# Recalculate so incidence2018_plus (which reads prevalence2018!Cxx/100)
# picks up the new figures.
$excel.Calculate()

# Update the view state on prevalence2018 to match where the edit focused.
$ws.Activate()
$ws.Range("C57:C122").Select()

# Update the view state on the "data" sheet (selection moved to the STEP
# reference figures used for the new prevalence numbers).
$dataWs = $wb.Worksheets.Item("data")
$dataWs.Activate()
$dataWs.Range("D5:D8").Select()

# Re-activate prevalence2018 as the visible tab, matching the original
# workbook (it was the tabSelected sheet).
$ws.Activate()
